$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update "想去人数" (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 272
$ws1.Range("F3").Value = 75
$ws1.Range("F5").Value = 7278
$ws1.Range("F6").Value = 5512
$ws1.Range("F11").Value = 243
$ws1.Range("F12").Value = 159

# Sheet "全部类型" (sheet4): update the same entries (duplicated data across sheets)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 272
$ws4.Range("F3").Value = 75
$ws4.Range("F5").Value = 7278
$ws4.Range("F6").Value = 5512
$ws4.Range("F11").Value = 243
$ws4.Range("F14").Value = 159
